# Update the "cryptos" price/volume table with freshly scraped values.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# D-column cells that look like plain decimal numbers (single '.') need a
# leading apostrophe so Excel keeps storing them as text (matching the sheet's
# existing convention) instead of silently converting them to a Double.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "38.783.59"
$ws.Range("E2").Value = "  +0.41%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "2.100.84"
$ws.Range("E3").Value = "  +0.40%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.10%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "'227.54"
$ws.Range("E5").Value = "  -0.51%  "

# --- Row 6 (XRP) ---
$ws.Range("E6").Value = "  +0.61%  "

# --- Row 7 (Solana) ---
$ws.Range("D7").Value = "'62.09"
$ws.Range("E7").Value = "  +2.64%  "

# --- Row 8 (USDC) ---
$ws.Range("E8").Value = "  +0.10%  "

# --- Row 9 (Cardano) ---
$ws.Range("E9").Value = "  +1.28%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("D10").Value = "'0.0842"
$ws.Range("E10").Value = "  +0.40%  "

# --- Row 11 (TRON) ---
$ws.Range("E11").Value = "  -0.39%  "

# --- Row 12 (Chainlink) ---
$ws.Range("D12").Value = "'15.81"
$ws.Range("E12").Value = "  +5.36%  "

# --- Row 13 (WrappedliquidstakedEther2.0) ---
$ws.Range("D13").Value = "2.412.64"
$ws.Range("E13").Value = "  +0.69%  "

# --- Row 14 (Avalanche) ---
$ws.Range("D14").Value = "'21.96"
$ws.Range("E14").Value = "  -1.48%  "

# --- Row 15 (Polygon) ---
$ws.Range("E15").Value = "  +1.21%  "

# --- Row 16 (Polkadot) ---
$ws.Range("D16").Value = "'5.48"
$ws.Range("E16").Value = "  +0.45%  "

# --- Row 17 (WrappedEther) ---
$ws.Range("D17").Value = "2.107.74"
$ws.Range("E17").Value = "  +1.40%  "

# --- Row 18 (WrappedBTC) ---
$ws.Range("D18").Value = "38.753.30"
$ws.Range("E18").Value = "  +0.53%  "

# --- Row 19 (Litecoin) ---
$ws.Range("D19").Value = "'71.69"
$ws.Range("E19").Value = "  +0.29%  "

# --- Row 20 (Uniswap) ---
$ws.Range("D20").Value = "'6.05"
$ws.Range("E20").Value = "  +0.14%  "

# --- Row 21 (ShibaInu) ---
$ws.Range("D21").Value = "0.0₃0843"
$ws.Range("E21").Value = "  +0.84%  "

# --- Row 22 (BitcoinCash) ---
$ws.Range("D22").Value = "'227.05"
$ws.Range("E22").Value = "  +0.72%  "

# --- Row 23 (Dai) ---
$ws.Range("E23").Value = "  +0.04%  "

# --- Row 24 (Toncoin) ---
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -4.21%  "

# --- Row 25 (PancakeSwap) ---
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  -1.21%  "

# --- Row 26 (Cosmos) ---
$ws.Range("D26").Value = "'9.67"
$ws.Range("E26").Value = "  +1.93%  "

# --- Row 27 (Monero) ---
$ws.Range("D27").Value = "'170.35"
$ws.Range("E27").Value = "  -0.26%  "

# --- Row 28 (Kaspa) ---
$ws.Range("E28").Value = "  +0.58%  "

# --- Row 29 (ImmutableX) ---
$ws.Range("D29").Value = "'1.42"
$ws.Range("E29").Value = "  +2.32%  "

# --- Row 30 (EthereumClassic) ---
$ws.Range("E30").Value = "  +0.92%  "

# --- Row 31 (WEMIXToken) ---
$ws.Range("D31").Value = "'2.54"
$ws.Range("E31").Value = "  +8.88%  "

# --- Row 32 (Stellar) ---
$ws.Range("E32").Value = "  +0.43%  "

# --- Row 33 (Filecoin) ---
$ws.Range("E33").Value = "  +1.72%  "

# --- Row 34 (InternetComputer(DFINITY)) ---
$ws.Range("D34").Value = "'4.80"
$ws.Range("E34").Value = "  +0.93%  "

# --- Row 35 (THORChain) ---
$ws.Range("D35").Value = "'7.12"
$ws.Range("E35").Value = "  +11.99%  "

# --- Row 36 (Hedera) ---
$ws.Range("E36").Value = "  +0.37%  "

# --- Row 37 (LidoDAOToken) ---
$ws.Range("D37").Value = "'2.36"
$ws.Range("E37").Value = "  -1.04%  "

# --- Row 38 (RenderToken) ---
$ws.Range("D38").Value = "'3.50"
$ws.Range("E38").Value = "  -0.89%  "

# --- Row 39 (BinanceUSD) ---
$ws.Range("E39").Value = "  +0.20%  "

# --- Row 40 (InjectiveProtocol) ---
$ws.Range("D40").Value = "'18.01"
$ws.Range("E40").Value = "  -2.65%  "

# --- Row 41 (VeChain) ---
$ws.Range("E41").Value = "  +3.54%  "

# --- Row 42 (Aave) ---
$ws.Range("D42").Value = "'101.62"
$ws.Range("E42").Value = "  +0.07%  "

# --- Row 43 (Maker) ---
$ws.Range("D43").Value = "1.525.43"
$ws.Range("E43").Value = "  -1.12%  "

# --- Row 44 (TrustWalletToken) ---
$ws.Range("E44").Value = "  +7.96%  "

# --- Row 45 (HuobiToken) ---
$ws.Range("E45").Value = "  -0.21%  "

# --- Rows 46 & 47: FraxShare and Cronos swapped ranking order, plus updated
#     price/volume figures ---
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'7.78"
$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0911"
$ws.Range("E47").Value = "  -1.68%  "

# --- Row 48 (ARBITRUM) ---
$ws.Range("E48").Value = "  +4.80%  "

# --- Row 49 (FTXToken) ---
$ws.Range("E49").Value = "  +0.94%  "

# --- Row 50 (MXToken) ---
$ws.Range("E50").Value = "  -1.18%  "

# --- Row 51 (RocketPoolETH) ---
$ws.Range("D51").Value = "2.299.90"
$ws.Range("E51").Value = "  +0.65%  "
